$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-and Post-Test Alerts")
$ws.Range("B5").Value = "1.25-2.25"
